# Updated cryptos list on Fri May 19 07:47:40 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table with
# the latest coinranking.com snapshot, and re-sorts three rank-adjacent pairs
# (Chainlink/Litecoin, TheSandbox/FraxShare, Quant/PaxDollar) whose relative
# order flipped, including their Coin name + Link cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.918.98"
$ws.Range("E2").Value = "  -1.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.810.11"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.29%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.57"
$ws.Range("E5").Value = "  -0.88%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.27%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4606"
$ws.Range("E7").Value = "  +2.99%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3703"
$ws.Range("E8").Value = "  -1.94%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07372"
$ws.Range("E9").Value = "  -0.49%  "

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8733"
$ws.Range("E10").Value = "  -0.83%  "

# Row 11 - Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.46"
$ws.Range("E11").Value = "  -2.02%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.794.73"
$ws.Range("E12").Value = "  -1.27%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.356"
$ws.Range("E13").Value = "  -1.51%  "

# Row 14 - Chainlink (was Litecoin)
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.526"
$ws.Range("E14").Value = "  -2.91%  "

# Row 15 - Litecoin (was Chainlink)
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.91"
$ws.Range("E15").Value = "  -1.34%  "

# Row 16 - TRON
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07043"
$ws.Range("E16").Value = "  -0.27%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  +0.30%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008694"
$ws.Range("E18").Value = "  -1.43%  "

# Row 19 - Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.18%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  -2.28%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "26.916.06"
$ws.Range("E21").Value = "  -1.47%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.325"
$ws.Range("E22").Value = "  -0.42%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -2.96%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.031.92"
$ws.Range("E24").Value = "  -0.75%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.895"
$ws.Range("E25").Value = "  -3.18%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.32"
$ws.Range("E26").Value = "  +0.23%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.36"
$ws.Range("E27").Value = "  -1.33%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.141"
$ws.Range("E28").Value = "  -6.08%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.300"
$ws.Range("E29").Value = "  -0.99%  "

# Row 30 - BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.78"
$ws.Range("E30").Value = "  -1.32%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08905"
$ws.Range("E31").Value = "  +0.33%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -4.88%  "

# Row 33 - ARBITRUM
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.154"
$ws.Range("E33").Value = "  -3.67%  "

# Row 34 - HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.918"
$ws.Range("E34").Value = "  -0.11%  "

# Row 35 - Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.445"
$ws.Range("E35").Value = "  -2.80%  "

# Row 36 - Frax
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.24%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.102"
$ws.Range("E37").Value = "  -0.61%  "

# Row 38 - VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01972"
$ws.Range("E38").Value = "  -0.02%  "

# Row 39 - Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05234"
$ws.Range("E39").Value = "  -0.60%  "

# Row 40 - RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.417"
$ws.Range("E40").Value = "  +2.73%  "

# Row 41 - MXToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.931"
$ws.Range("E41").Value = "  +1.91%  "

# Row 42 - FraxShare (was TheSandbox)
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.178"
$ws.Range("E42").Value = "  -1.65%  "

# Row 43 - TheSandbox (was FraxShare)
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5287"
$ws.Range("E43").Value = "  -0.30%  "

# Row 44 - Algorand
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1659"
$ws.Range("E44").Value = "  -2.61%  "

# Row 45 - Aptos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.502"
$ws.Range("E45").Value = "  -1.58%  "

# Row 46 - Decentraland
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4970"
$ws.Range("E46").Value = "  -1.74%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.30"
$ws.Range("E47").Value = "  -2.78%  "

# Row 48 - PaxDollar (was Quant)
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.30%  "

# Row 49 - Quant (was PaxDollar)
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.98"
$ws.Range("E49").Value = "  -1.49%  "

# Row 50 - NEARProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.669"
$ws.Range("E50").Value = "  -1.01%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06296"
$ws.Range("E51").Value = "  -1.41%  "
